$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column J: header label "18_22" and repeat the 2022 (column I) values.
$ws.Range("J7").Value2 = "18_22"
$ws.Range("J8").Value2 = $ws.Range("I8").Value2
$ws.Range("J9").Value2 = $ws.Range("I9").Value2
$ws.Range("J10").Value2 = $ws.Range("I10").Value2
$ws.Range("J11").Value2 = $ws.Range("I11").Value2
$ws.Range("J12").Value2 = $ws.Range("I12").Value2
$ws.Range("J13").Value2 = $ws.Range("I13").Value2
$ws.Range("J14").Value2 = $ws.Range("I14").Value2

$ws.Range("H18").Select()
